$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Move the <w:lastRenderedPageBreak/> marker that currently sits at the
#    start of the "Linux version 3.9.0-xilinx ..." run so that it instead
#    sits at the start of the "This one seemed important. Sets up our
#    clock." run further down the same section.
# ---------------------------------------------------------------------------

# 1a. Remove it from the "Linux version ..." run (delete + reinsert the run
#     text without the page-break marker, keeping the bold formatting).
$r = $d.Content
$found = $r.Find.Execute("Linux version 3.9.0-xilinx (kfisch13@linux-4.ece.iastate.edu) (", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find Linux version run" }
$r.Delete()
$collapsed = $d.Range($r.Start, $r.Start)
$collapsed.InsertXML('<w:p><w:r><w:rPr><w:b/></w:rPr><w:t>Linux version 3.9.0-xilinx (kfisch13@linux-4.ece.iastate.edu) (</w:t></w:r></w:p>')

# 1b. Add it to the "This one seemed important..." run.
$r = $d.Content
$found = $r.Find.Execute("This one seemed important. Sets up our clock.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'This one seemed important' run" }
$r.Delete()
$collapsed = $d.Range($r.Start, $r.Start)
$collapsed.InsertXML('<w:p><w:r><w:lastRenderedPageBreak/><w:t>This one seemed important. Sets up our clock.</w:t></w:r></w:p>')

# ---------------------------------------------------------------------------
# 2) The "Changes to launcher_driver.c" Heading1 now starts a new page, so
#    Word stamps a fresh <w:lastRenderedPageBreak/> on its first run.
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Changes to ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'Changes to ' run" }
$r.Delete()
$collapsed = $d.Range($r.Start, $r.Start)
$collapsed.InsertXML('<w:p><w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">Changes to </w:t></w:r></w:p>')

# ---------------------------------------------------------------------------
# 3) Rewrite the last bullet of "Changes to launcher_driver.c":
#      "Replaced " + bookmark(_GoBack) + "the IDs in line 33 ... launcher_commands.h"
#    becomes
#      "Replaced the IDs in line 33 ... launcher_commands.h"   (bookmark removed)
#    and add the new "How launcher_fire.c Works" section (with the _GoBack
#    bookmark moved onto its heading) right after it.
# ---------------------------------------------------------------------------
$r = $d.Content
$found = $r.Find.Execute("Replaced*launcher_commands.h", $true, $false, $true, $false, $false, $true, 1, $false, "", 0)
if (-not $found) { throw "Could not find 'Replaced ... launcher_commands.h' bullet" }
# include the paragraph mark so the whole bullet paragraph disappears
$extended = $d.Range($r.Start, $r.End + 1)
$extended.Delete()

$collapsed = $d.Range($extended.Start, $extended.Start)
$frag = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Replaced the IDs in line 33 to their appropriate defines in </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>launcher_commands.h</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p/><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:r><w:t xml:space="preserve">How </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>launcher_fire.c</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> Works</w:t></w:r><w:bookmarkStart w:id="100" w:name="_GoBack"/><w:bookmarkEnd w:id="100"/></w:p><w:p><w:r><w:t>This file has a main function and a helper function. The helper function has some debug information about the launch and it has the actual driver call in it. The main function gets a reference to the driver, calls the helper function to start and stop firing, and closes the reference to the driver.</w:t></w:r></w:p><w:p/><w:p/>'
$collapsed.InsertXML($frag)

Write-Output "edit.ps1 completed"
